$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates: force text format to avoid Excel auto-converting
# numeric-looking strings (e.g. "4.26") into actual numbers.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "25.829.71"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.636.47"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "215.95"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.0639"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "19.64"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "4.26"
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.632.00"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "1.862.45"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.553"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0₃0773"
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "63.39"
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "25.846.66"
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "193.80"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.79"
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "139.89"
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "15.64"
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.24"
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0488"
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.898"
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.552"
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.109.29"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "99.67"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.804"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.42"
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "7.68"
$c.Style = "Normal"

# Other column (B/C/E) updates: plain text assignment is safe.
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -4.32%  "
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("E47").Value = "  +9.78%  "
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("E51").Value = "  +0.53%  "
